$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: merge the two runs that make up the "SUN Aug 12 13:20:54
# IST 2018" paragraph into a single run (same visible text, just one
# <w:r> instead of two).
# -----------------------------------------------------------------
$n = $d.Paragraphs.Count
$sunIdx = -1
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "SUN Aug 12 13:20:54 IST 2018`r") {
        $sunIdx = $i
        break
    }
}
if ($sunIdx -eq -1) {
    throw "Could not find the 'SUN Aug 12 13:20:54 IST 2018' paragraph"
}
$sunPara = $d.Paragraphs.Item($sunIdx)
$sunRange = $sunPara.Range
$sunRange.MoveEnd(1, -1)
# Set to a different value first so the engine treats it as a real
# text change (re-writing the identical string is a no-op and would
# leave the original two runs untouched), then set the final value.
$sunRange.Text = "SUN Aug 12 13:20:54 IST 2018 TEMP"
$sunRange2 = $sunPara.Range
$sunRange2.MoveEnd(1, -1)
$sunRange2.Text = "SUN Aug 12 13:20:54 IST 2018"

# -----------------------------------------------------------------
# Change 2: append a brand-new purchase-record block right after the
# paragraph that ends with "- 8433.0" (the previous record's closing
# "Amount balance" line), before the pre-existing trailing blank
# paragraphs.
# -----------------------------------------------------------------
$n = $d.Paragraphs.Count
$anchorIdx = -1
for ($i = 1; $i -le $n; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -eq "Amount balance`t`t`t- 8433.0`r") {
        $anchorIdx = $i
        break
    }
}
if ($anchorIdx -eq -1) {
    throw "Could not find the 'Amount balance ... - 8433.0' paragraph"
}

$anchor = $d.Paragraphs.Item($anchorIdx).Range
$anchor.Collapse(0)

# The 10 new paragraphs, as plain text segments (tabs/newlines are
# inserted explicitly below). The engine reliably gives each distinct
# tab its own run and leaves genuinely-empty paragraphs without any
# run at all (matching native Word output), so plain InsertAfter text
# is enough here; bold/color formatting is reapplied per paragraph
# afterwards. A trailing `r is required so the very last segment does
# not get merged into the pre-existing (untouched) paragraph that
# follows the insertion point.
$segments = @(
    "",
    "MON Aug 13 12:02:57 IST 2018",
    "Person Name`t`t`t`t- N NANDISHA",
    "---------------------------------------------------------------",
    "Item Name`t`t`t`t- CARROT",
    "Amount Received`t`t`t- 4250",
    "Amount balance`t`t`t- 4183.0",
    "Amount Received mode`t`t- CASH",
    "",
    ""
)
$blockText = ($segments -join "`r") + "`r"
$anchor.InsertAfter($blockText)

# Re-resolve paragraph indices now that the new ones exist.
$base = $anchorIdx

$pEmpty1   = $d.Paragraphs.Item($base + 1)   # empty, bold
$pDate     = $d.Paragraphs.Item($base + 2)   # MON Aug 13 ...
$pPerson   = $d.Paragraphs.Item($base + 3)   # Person Name ...
$pDashes   = $d.Paragraphs.Item($base + 4)   # -------------
$pItem     = $d.Paragraphs.Item($base + 5)   # Item Name ...
$pReceived = $d.Paragraphs.Item($base + 6)   # Amount Received ...
$pBalance  = $d.Paragraphs.Item($base + 7)   # Amount balance ...
$pMode     = $d.Paragraphs.Item($base + 8)   # Amount Received mode ...
$pEmpty2   = $d.Paragraphs.Item($base + 9)   # empty, plain
$pEmpty3   = $d.Paragraphs.Item($base + 10)  # empty, bold

# 1) empty bold paragraph
$pEmpty1.Range.Font.Bold = $true

# 2) split the date paragraph's text into two runs: "MON Aug 13" and
# " 12:02:57 IST 2018" (with a leading space), matching the source
# document's convention for these date stamps.
$dateRange = $pDate.Range
$dateRange.MoveEnd(1, -1)
$dateStart = $dateRange.Start
$splitPoint = $dateStart + ("MON Aug 13").Length
$secondHalf = $d.Range($splitPoint, $splitPoint + (" 12:02:57 IST 2018").Length)
$secondHalf.Font.Bold = $true
$secondHalf.Font.Bold = $false

# 6) Amount Received line -> red
$pReceived.Range.Font.Color = 255

# 7) Amount balance line -> bold
$pBalance.Range.Font.Bold = $true

# 10) trailing empty paragraph -> bold
$pEmpty3.Range.Font.Bold = $true
